$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.815.33'
$ws.Range("E2").Value = '  -0.60%  '
$ws.Range("D3").Value = '2.235.03'
$ws.Range("E3").Value = '  -1.84%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.04'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.24%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '98.20'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -7.22%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.571'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.17%  '
$ws.Range("E8").Value = '  +0.18%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.529'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -7.58%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.75'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -8.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0819'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.57%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.33'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -7.28%  '
$ws.Range("E13").Value = '  -2.75%  '
$ws.Range("D14").Value = '2.574.66'
$ws.Range("E14").Value = '  -1.93%  '
$ws.Range("D15").Value = '2.235.98'
$ws.Range("E15").Value = '  -2.08%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.835'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.27%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.84'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.44%  '
$ws.Range("D18").Value = '43.723.01'
$ws.Range("E18").Value = '  -0.65%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.06'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -7.41%  '
$ws.Range("D20").Value = '0.0₃0967'
$ws.Range("E20").Value = '  -3.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.28'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.14%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.74'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.67%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.23'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.65%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.97'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -7.45%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.01'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -8.65%  '
$ws.Range("E26").Value = '  +0.26%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.06'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.14%  '
$ws.Range("E28").Value = '  -2.33%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '36.40'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -6.33%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.94'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -8.82%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '19.99'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.11%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '154.98'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.10%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0827'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.37%  '
$ws.Range("E34").Value = '  +1.45%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.64'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.30%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.90'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -8.61%  '
$ws.Range("E37").Value = '  -5.86%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '15.54'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.26%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.50'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -12.06%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.98'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -11.27%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0306'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.28%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.01'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.27%  '
$ws.Range("D44").Value = '1.696.10'
$ws.Range("E44").Value = '  -3.96%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '82.52'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.71%  '
$ws.Range("E46").Value = '  -6.85%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.14'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.75%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '101.20'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.11%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '70.98'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.60%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '55.96'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -6.11%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.59'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.20%  '
